$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Phase 1 Pre CPP")

# Copy the header-row formatting (bold font, thin border, center/top alignment)
# from A1 onto the new date cell A2, then set its value and apply the
# mm-dd-yyyy number format (this creates numFmtId 164 and a new cellXfs entry).
$ws.Range("A1").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("A2").Value = 30414
$ws.Range("A2").NumberFormat = "mm-dd-yyyy"

$ws.Range("B2").Value = 0.03033333333333333
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 261.8439554330322
$ws.Range("E2").Value = 0.009791527531336826

$excel.CutCopyMode = $false
